# working on logging bug
#
# Diary update: two new work-log entries were appended to the bottom of the
# "Arbeitszeit" sheet (rows 49 & 50), following the exact same layout as
# every other entry (Datum/Zeit/Einheit/Tätigkeit/Kommentar in columns
# E-I). The running totals in column B recalculate automatically since
# they're formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 49: 05.09.2019, 4 Stunden Programmieren ---
$ws.Cells.Item(49, 5).Value = 43713                         # E49 Datum
$ws.Cells.Item(49, 6).Value = 4                             # F49 Zeit
$ws.Cells.Item(49, 7).Value = "Stunden"                     # G49 Einheit
$ws.Cells.Item(49, 8).Value = "Programmieren"                # H49 Tätigkeit
$ws.Cells.Item(49, 9).Value = "Änderungen im Basis Code"    # I49 Kommentar

# --- New row 50: 06.09.2019, 3 Stunden Programmieren ---
$ws.Cells.Item(50, 5).Value = 43714                         # E50 Datum
$ws.Cells.Item(50, 6).Value = 3                             # F50 Zeit
$ws.Cells.Item(50, 7).Value = "Stunden"                     # G50 Einheit
$ws.Cells.Item(50, 8).Value = "Programmieren"                # H50 Tätigkeit
$ws.Cells.Item(50, 9).Value = "Weiterführende Anpasseungen" # I50 Kommentar

# Give the two new date cells the same number format as the rest of the
# "Datum" column (reuses the existing date style instead of minting a new
# one).
$ws.Range("E48").Copy()
$ws.Range("E49:E50").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Leave the sheet scrolled/selected the way the author left it ---
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("B11").Select()
